# Updates the cryptos worksheet with newly scraped price/volume data.
# Generated to match "Updated cryptos list" GitHub Actions commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values that look like plain numbers (single decimal point)
# are prefixed with a leading apostrophe so Excel stores them as text,
# matching the source data which is always a text string (e.g. "96.09"
# rather than the number 96.09). Values that already contain multiple dots
# (e.g. "43.725.75") are naturally stored as text and need no prefix.

$updates = @(
    @{Cell="D2";  Value="43.725.75"},
    @{Cell="E2";  Value="  +0.28%  "},

    @{Cell="D3";  Value="2.290.12"},
    @{Cell="E3";  Value="  -1.43%  "},

    @{Cell="E4";  Value="  -0.09%  "},

    @{Cell="D5";  Value="'96.09"},
    @{Cell="E5";  Value="  +2.43%  "},

    @{Cell="D6";  Value="'269.53"},
    @{Cell="E6";  Value="  +0.30%  "},

    @{Cell="D7";  Value="'0.624"},
    @{Cell="E7";  Value="  +0.52%  "},

    @{Cell="E8";  Value="  -0.15%  "},

    @{Cell="E9";  Value="  -1.39%  "},

    @{Cell="D10"; Value="'45.59"},
    @{Cell="E10"; Value="  +2.33%  "},

    @{Cell="D11"; Value="'0.0932"},
    @{Cell="E11"; Value="  -0.87%  "},

    @{Cell="E12"; Value="  -1.03%  "},

    @{Cell="D13"; Value="'0.106"},
    @{Cell="E13"; Value="  +1.34%  "},

    @{Cell="D14"; Value="'15.68"},
    @{Cell="E14"; Value="  +1.61%  "},

    @{Cell="D15"; Value="2.634.65"},
    @{Cell="E15"; Value="  -1.19%  "},

    # Rows 16 and 17 swap (Polygon <-> WrappedEther) along with new data
    @{Cell="B16"; Value="WrappedEther"},
    @{Cell="C16"; Value="https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"},
    @{Cell="D16"; Value="2.431.20"},
    @{Cell="E16"; Value="  +4.97%  "},

    @{Cell="B17"; Value="Polygon"},
    @{Cell="C17"; Value="https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"},
    @{Cell="D17"; Value="'0.849"},
    @{Cell="E17"; Value="  -1.12%  "},

    @{Cell="D18"; Value="43.673.53"},
    @{Cell="E18"; Value="  +0.23%  "},

    @{Cell="E19"; Value="  +3.35%  "},

    @{Cell="E20"; Value="  -2.77%  "},

    @{Cell="D21"; Value="'72.11"},

    @{Cell="D22"; Value="'2.48"},
    @{Cell="E22"; Value="  +9.86%  "},

    @{Cell="D23"; Value="'232.82"},
    @{Cell="E23"; Value="  -1.82%  "},

    @{Cell="D24"; Value="'9.09"},
    @{Cell="E24"; Value="  -5.18%  "},

    @{Cell="D25"; Value="'2.70"},
    @{Cell="E25"; Value="  +8.01%  "},

    @{Cell="E26"; Value="  -0.17%  "},

    @{Cell="D27"; Value="'11.31"},
    @{Cell="E27"; Value="  +0.03%  "},

    @{Cell="E28"; Value="  -0.61%  "},

    # Rows 29 and 30 swap (Toncoin <-> InjectiveProtocol) along with new data
    @{Cell="B29"; Value="InjectiveProtocol"},
    @{Cell="C29"; Value="https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"},
    @{Cell="D29"; Value="'38.82"},
    @{Cell="E29"; Value="  +1.00%  "},

    @{Cell="B30"; Value="Toncoin"},
    @{Cell="C30"; Value="https://coinranking.com/coin/67YlI0K1b+toncoin-ton"},
    @{Cell="D30"; Value="'2.23"},
    @{Cell="E30"; Value="  -2.54%  "},

    @{Cell="D31"; Value="'174.85"},
    @{Cell="E31"; Value="  +1.82%  "},

    @{Cell="D32"; Value="'22.06"},
    @{Cell="E32"; Value="  -3.00%  "},

    @{Cell="E33"; Value="  +0.25%  "},

    @{Cell="E35"; Value="  +0.33%  "},

    @{Cell="D36"; Value="'4.58"},
    @{Cell="E36"; Value="  +5.74%  "},

    @{Cell="E37"; Value="  -0.48%  "},

    @{Cell="E38"; Value="  -1.27%  "},

    @{Cell="E39"; Value="  +5.10%  "},

    @{Cell="E40"; Value="  +1.03%  "},

    @{Cell="E41"; Value="  -2.30%  "},

    @{Cell="D42"; Value="'12.32"},
    @{Cell="E42"; Value="  +2.40%  "},

    @{Cell="E43"; Value="  -1.30%  "},

    @{Cell="D44"; Value="'64.51"},
    @{Cell="E44"; Value="  +4.91%  "},

    # Rows 45 and 46 swap (FraxShare <-> THORChain) along with new data
    @{Cell="B45"; Value="THORChain"},
    @{Cell="C45"; Value="https://coinranking.com/coin/ybmU-kKU+thorchain-rune"},
    @{Cell="D45"; Value="'5.20"},
    @{Cell="E45"; Value="  -2.75%  "},

    @{Cell="B46"; Value="FraxShare"},
    @{Cell="C46"; Value="https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"},
    @{Cell="D46"; Value="'8.72"},
    @{Cell="E46"; Value="  -3.77%  "},

    @{Cell="E47"; Value="  -0.19%  "},

    @{Cell="E48"; Value="  -0.34%  "},

    @{Cell="E49"; Value="  -2.81%  "},

    @{Cell="D50"; Value="'1.53"},
    @{Cell="E50"; Value="  +12.33%  "},

    @{Cell="D51"; Value="'0.430"},
    @{Cell="E51"; Value="  +3.76%  "}
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
